$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple address corrections ---
$ws.Range("E2").Value = "高松市観光通2-9-15"
$ws.Range("E21").Value = "高松市仏生山町甲2518-16"
$ws.Range("E42").Value = "高松市東山崎町33‐2"

# --- Remove the row for "塩江分団第２部岩部屯所" (old row 64) ---
# This shifts rows 65-97 up to 64-96.
$ws.Rows(64).Delete()

# --- Row 62 ("塩江分団第１部長野屯所") becomes a new record "塩江分団第１部車庫" ---
# Target row has no A (number) and no H (note) cell at all, only B-G populated.
$ws.Range("A62").ClearContents()
$ws.Range("B62").Value = "34.19596225788533"
$ws.Range("C62").Value = "134.04199847291943"
$ws.Range("D62").Value = "高松市消防団塩江分団第１部車庫"
$ws.Range("E62").Value = "高松市塩江町安原下第２号1645-1"
$ws.Range("H62").ClearContents()

# --- Row 95 (after the deletion above, was old row 96 "消防山田倉庫") becomes ---
# --- a new record "女木分団西浦消防屯所倉庫", again without A and H cells. ---
$ws.Range("A95").ClearContents()
$ws.Range("B95").Value = "34.39506251264558"
$ws.Range("C95").Value = "134.04450044702173"
$ws.Range("D95").Value = "高松市消防団女木分団西浦消防屯所倉庫"
$ws.Range("E95").Value = "高松市女木町西浦1906地先"
$ws.Range("H95").ClearContents()

# --- Row 96 (old row 97 "香川分団第３部川東下屯所") keeps its data but the ---
# --- street address house number changes from 765-6 to 765-7. ---
$ws.Range("E96").Value = "高松市香川町川東下765-7"

"script done"
